# File3.xlsx — populate Аркуш1 (A1:C2) with the 2x3 number grid added in the
# commit, and leave the selection on E5 (matches the saved <selection>).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 1
$ws.Cells.Item(1, 2).Value = 2
$ws.Cells.Item(1, 3).Value = 3
$ws.Cells.Item(2, 1).Value = 4
$ws.Cells.Item(2, 2).Value = 5
$ws.Cells.Item(2, 3).Value = 6

[void]$ws.Range("E5").Select()
